# Applies the changes described by the diff:
#  1. Swap the data values of columns A (year) and B (month) for rows 2-13,
#     so that A becomes month and B becomes year (headers stay as-is).
#  2. Add new header columns H1:P1 with grade/frequency labels.
#  3. Add new data values for H2:P13 (same constant values on every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap A/B data values for rows 2 through 13 ---
for ($r = 2; $r -le 13; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $bVal
    $ws.Cells.Item($r, 2).Value2 = $aVal
}

# --- Step 2: add new headers in H1:P1 ---
$ws.Range("H1").Value2 = "grade_total"
$ws.Range("I1").Value2 = "grade_distance"
$ws.Range("J1").Value2 = "grade_visitation"
$ws.Range("K1").Value2 = "grade_encounters"
$ws.Range("L1").Value2 = "NEVER"
$ws.Range("M1").Value2 = "RARELY"
$ws.Range("N1").Value2 = "SOMETIMES"
$ws.Range("O1").Value2 = "FREQUENTLY"
$ws.Range("P1").Value2 = "ALWAYS"

# --- Step 3: add new data values for H2:P13 (constant across rows) ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value2 = 0        # H - grade_total
    $ws.Cells.Item($r, 9).Value2 = 1        # I - grade_distance
    $ws.Cells.Item($r, 10).Value2 = 0       # J - grade_visitation
    $ws.Cells.Item($r, 11).Value2 = 0       # K - grade_encounters
    $ws.Cells.Item($r, 12).Value2 = 1.023   # L - NEVER
    $ws.Cells.Item($r, 13).Value2 = 1.021   # M - RARELY
    $ws.Cells.Item($r, 14).Value2 = 1.072   # N - SOMETIMES
    $ws.Cells.Item($r, 15).Value2 = 1.162   # O - FREQUENTLY
    $ws.Cells.Item($r, 16).Value2 = 1.722   # P - ALWAYS
}
